$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old row 13 held a stray value ("210064 - Eduardo Rezende Triboni") in
# columns B/C with no label in column A. That row is removed, which shifts
# every following row up by one position and drops the dimension by a row.
$ws.Rows(13).Delete()

# A handful of the shifted rows need their B/C values corrected so the
# label in column A lines up with the right value again.
$ws.Range("B10:C10").Value = "210064 - Eduardo Rezende Triboni"
$ws.Range("B13:C13").Value = "Semestral"

# Row 15 needs to show a literal date-like string ("01/01/2020"); assigning
# it directly would make Excel reinterpret it as a date serial number, so
# instead copy the existing text value from B8 (which already holds that
# same string) to keep it stored as text.
$ws.Range("B8").Copy() | Out-Null
$ws.Range("B15").PasteSpecial(-4163) | Out-Null
$ws.Range("B8").Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4163) | Out-Null

$ws.Range("B18:C18").Value = "210064 - Eduardo Rezende Triboni"
$ws.Range("B19:C19").Value = "Duas provas teóricas e ao longo do semestre letivoAos alunos que tiverem freqüência mínima de 70% e média final menor que 5,0 e igual ou maior que 3,0, será dada recuperação com uma avaliação escrita. A média dessa avaliação somada com a média anterior das P1 e P2, se superior a cinco (5,0), levará a aprovação do aluno."
$ws.Range("B20:C20").Value = "A média final (M) será calculada pela expressão: M = (P1 + P2)/2"
$ws.Range("B21:C21").Value = "Aos alunos que tiverem freqüência mínima de 70% e média final menor que 5,0 e igual ou maior que 3,0, será dada recuperação com uma avaliação escrita. A média dessa avaliação somada com a média anterior das P1 e P2, se superior a cinco (5,0), levará a aprovação do aluno."

$excel.CutCopyMode = 0
